$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K, shifting existing K:R columns to L:S.
# This preserves all existing formatting/values/styles of the shifted cells.
$ws.Columns("K:K").Insert()

# New column header + data: "VAT Partner" info column (header row only;
# the data row K2 is left blank, matching the template row).
$ws.Range("K1").Value = "VAT Partner"

# Match the neighbouring column's width for the newly inserted column.
$ws.Columns("K:K").ColumnWidth = $ws.Columns("J:J").ColumnWidth

# The worksheet's hidden _FilterDatabase defined name must grow to cover
# the new column (was $A$1:$R$1, now $A$1:$S$1).
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$S`$1"
    }
}

# Match the saved selection state (active cell K2) from the source file.
[void]$ws.Range("K2").Select()
